$wb = $excel.ActiveWorkbook

# --- Rename the original sheet, add the two new ones (in final left-to-right order) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "ArcFace"

$ws2 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws2.Name = "VGGFace"

$ws3 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3.Name = "FaceNet512"

# ---------------------------------------------------------------------------
# ArcFace
# ---------------------------------------------------------------------------
$ws1.Range("A1").Value = "Metric"
$ws1.Range("B1").Value = "Value (Weighted)"
$ws1.Range("C1").Value = "Value (Micro)"
$ws1.Range("D1").Value = "Value(Macro)"

$ws1.Range("A2").Value = "Accuracy"
$ws1.Range("B2").Value = 0.94514100000000001
$ws1.Range("C2").Value = 0.94514100000000001
$ws1.Range("D2").Value = 0.94514100000000001

$ws1.Range("A3").Value = "Precision"
$ws1.Range("B3").Value = 0.95150199999999996
$ws1.Range("C3").Value = 0.94514100000000001
$ws1.Range("D3").Value = 0.63668400000000003

$ws1.Range("A4").Value = "Recall"
$ws1.Range("B4").Value = 0.94514100000000001
$ws1.Range("C4").Value = 0.94514100000000001
$ws1.Range("D4").Value = 0.62698399999999999

$ws1.Range("A5").Value = "F1-Score"
$ws1.Range("B5").Value = 0.94544399999999995
$ws1.Range("C5").Value = 0.94514100000000001
$ws1.Range("D5").Value = 0.62987199999999999

$ws1.Columns.Item(2).ColumnWidth = 15.26953125
$ws1.Columns.Item(3).ColumnWidth = 12
$ws1.Columns.Item(4).ColumnWidth = 12.08984375

# ---------------------------------------------------------------------------
# VGGFace
# ---------------------------------------------------------------------------
$ws2.Range("A1").Value = "Metric"
$ws2.Range("B1").Value = "Value (Weighted)"
$ws2.Range("C1").Value = "Value (Micro)"
$ws2.Range("D1").Value = "Value(Macro)"

$ws2.Range("A2").Value = "Accuracy"
$ws2.Range("B2").Value = 0.99216300000000002
$ws2.Range("C2").Value = 0.99216300000000002
$ws2.Range("D2").Value = 0.99216300000000002

$ws2.Range("A3").Value = "Precision"
$ws2.Range("B3").Value = 0.99380199999999996
$ws2.Range("C3").Value = 0.99216300000000002
$ws2.Range("D3").Value = 0.66283499999999995

$ws2.Range("A4").Value = "Recall"
$ws2.Range("B4").Value = 0.99216300000000002
$ws2.Range("C4").Value = 0.99216300000000002
$ws2.Range("D4").Value = 0.66099799999999997

$ws2.Range("A5").Value = "F1-Score"
$ws2.Range("B5").Value = 0.99293100000000001
$ws2.Range("C5").Value = 0.99216300000000002
$ws2.Range("D5").Value = 0.66188100000000005

$ws2.Columns.Item(2).ColumnWidth = 15.26953125
$ws2.Columns.Item(3).ColumnWidth = 12
$ws2.Columns.Item(4).ColumnWidth = 12.08984375

# ---------------------------------------------------------------------------
# FaceNet512
# ---------------------------------------------------------------------------
$ws3.Range("A1").Value = "Metric"
$ws3.Range("B1").Value = "Value (Weighted)"
$ws3.Range("C1").Value = "Value (Micro)"
$ws3.Range("D1").Value = "Value(Macro)"

$ws3.Range("A2").Value = "Accuracy"
$ws3.Range("B2").Value = 0.996865
$ws3.Range("C2").Value = 0.996865
$ws3.Range("D2").Value = 0.996865

$ws3.Range("A3").Value = "Precision"
$ws3.Range("B3").Value = 0.99843700000000002
$ws3.Range("C3").Value = 0.996865
$ws3.Range("D3").Value = 0.66569999999999996

$ws3.Range("A4").Value = "Recall"
$ws3.Range("B4").Value = 0.996865
$ws3.Range("C4").Value = 0.996865
$ws3.Range("D4").Value = 0.66439899999999996

$ws3.Range("A5").Value = "F1-Score"
$ws3.Range("B5").Value = 0.997645
$ws3.Range("C5").Value = 0.996865
$ws3.Range("D5").Value = 0.665045

$ws3.Columns.Item(2).ColumnWidth = 15.26953125
$ws3.Columns.Item(3).ColumnWidth = 12
$ws3.Columns.Item(4).ColumnWidth = 12.08984375

# --- Selections on each sheet (match final saved state) ---
$ws1.Range("D7").Select()
$ws2.Range("E4").Select()
$ws3.Range("E3").Select()

# --- FaceNet512 is the active/visible tab when the workbook is saved ---
$ws3.Select()
